$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" column values look numeric (e.g. "7.46", "0.100") but the
# source workbook stores every Coin/Link/Price/Volume cell as text. Force a
# Text number format before assigning these so Excel keeps them as text
# instead of silently converting to a number (which would also drop
# significant trailing zeros, e.g. "0.100" -> 0.1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.355.50"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.591.97"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "210.23"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.813.96"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.08"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.565.89"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "64.70"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "26.356.17"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "7.46"
$ws.Range("E19").Value = "  +4.36%  "
$ws.Range("D20").Value = "212.21"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "8.92"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "144.83"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Value = "1.297.42"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").Value = "0.612"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  -11.20%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D42").Value = "5.62"
$ws.Range("E42").Value = "  +2.97%  "
$ws.Range("D43").Value = "62.77"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "0.762"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "1.727.15"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  -1.18%  "
